# Updated cryptos list on Tue Jun 11 13:47:17 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    # Force text interpretation so numeric-looking strings (e.g. "0.482",
    # "67.016.91") are preserved verbatim instead of being parsed as numbers.
    $rng.NumberFormat = "@"
    $rng.Value = $value
    # Restore the default/original cell style so we don't leave a stray
    # "Text" number format applied to the cell's style record.
    $rng.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "66.994.44"
Set-TextValue "E2" "  -3.58%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.523.40"
Set-TextValue "E3" "  -4.21%  "

# Row 4 - TetherUSD
Set-TextValue "E4" "  +0.06%  "

# Row 5 - BNB
Set-TextValue "D5" "607.25"
Set-TextValue "E5" "  -6.03%  "

# Row 6 - Solana
Set-TextValue "D6" "151.84"
Set-TextValue "E6" "  -4.51%  "

# Row 7 - LidoStakedEther
Set-TextValue "D7" "3.521.09"
Set-TextValue "E7" "  -4.24%  "

# Row 8 - USDC
Set-TextValue "E8" "  +0.22%  "

# Row 9 - XRP
Set-TextValue "D9" "0.482"
Set-TextValue "E9" "  -3.29%  "

# Row 10 - Dogecoin
Set-TextValue "D10" "0.139"
Set-TextValue "E10" "  -3.96%  "

# Row 11 - Toncoin
Set-TextValue "D11" "6.78"
Set-TextValue "E11" "  -4.66%  "

# Row 12 - Cardano
Set-TextValue "D12" "0.425"
Set-TextValue "E12" "  -3.62%  "

# Row 13 - ShibaInu
Set-TextValue "D13" "0.0000219"
Set-TextValue "E13" "  -4.92%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-TextValue "D14" "4.127.62"
Set-TextValue "E14" "  -3.96%  "

# Row 15 - Avalanche
Set-TextValue "D15" "31.51"
Set-TextValue "E15" "  -3.33%  "

# Row 16 - WrappedEther
Set-TextValue "D16" "3.539.79"
Set-TextValue "E16" "  -3.29%  "

# Row 17 - WrappedBTC
Set-TextValue "D17" "67.016.91"
Set-TextValue "E17" "  -3.53%  "

# Row 18 - TRON
Set-TextValue "E18" "  +0.91%  "

# Row 19 - Polkadot
Set-TextValue "D19" "6.27"
Set-TextValue "E19" "  -3.06%  "

# Row 20 - Chainlink
Set-TextValue "D20" "15.33"
Set-TextValue "E20" "  -3.57%  "

# Row 21 - BitcoinCash
Set-TextValue "D21" "442.03"
Set-TextValue "E21" "  -5.18%  "

# Row 22 - Uniswap
Set-TextValue "D22" "9.19"
Set-TextValue "E22" "  -8.31%  "

# Row 23 - Polygon
Set-TextValue "D23" "0.624"
Set-TextValue "E23" "  -3.52%  "

# Row 24 - Litecoin
Set-TextValue "D24" "77.60"
Set-TextValue "E24" "  -2.27%  "

# Row 25 - WrappedeETH
Set-TextValue "D25" "3.673.00"
Set-TextValue "E25" "  -3.91%  "

# Row 26 - Dai
Set-TextValue "E26" "  -0.05%  "

# Row 27 - PEPE
Set-TextValue "D27" "0.0000120"
Set-TextValue "E27" "  -3.72%  "

# Row 28 - InternetComputer(DFINITY)
Set-TextValue "D28" "10.11"
Set-TextValue "E28" "  -5.88%  "

# Row 29 - RenderToken
Set-TextValue "D29" "8.15"
Set-TextValue "E29" "  -9.59%  "

# Row 30 - PancakeSwap
Set-TextValue "D30" "2.52"
Set-TextValue "E30" "  -4.12%  "

# Row 31 - Fetch.AI
Set-TextValue "D31" "1.66"
Set-TextValue "E31" "  -2.10%  "

# Row 32 - Binance-PegBSC-USD
Set-TextValue "E32" "  +0.05%  "

# Row 33 - EthereumClassic
Set-TextValue "D33" "25.54"
Set-TextValue "E33" "  -5.22%  "

# Row 34 - Kaspa
Set-TextValue "D34" "0.157"
Set-TextValue "E34" "  -3.29%  "

# Row 35 - was ImmutableX, now RenzoRestakedETH
Set-TextValue "B35" "RenzoRestakedETH"
Set-TextValue "C35" "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
Set-TextValue "D35" "3.525.07"
Set-TextValue "E35" "  -3.94%  "

# Row 36 - NEARProtocol
Set-TextValue "D36" "6.11"
Set-TextValue "E36" "  -4.96%  "

# Row 37 - was RenzoRestakedETH, now ImmutableX
Set-TextValue "B37" "ImmutableX"
Set-TextValue "C37" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D37" "1.85"
Set-TextValue "E37" "  -7.48%  "

# Row 38 - Aptos
Set-TextValue "D38" "7.95"
Set-TextValue "E38" "  -5.60%  "

# Row 39 - USDe
Set-TextValue "E39" "  +0.09%  "

# Row 40 - FirstDigitalUSD
Set-TextValue "E40" "  +0.08%  "

# Row 41 - Monero
Set-TextValue "D41" "175.11"
Set-TextValue "E41" "  -2.05%  "

# Row 42 - Stacks
Set-TextValue "D42" "2.13"
Set-TextValue "E42" "  -4.08%  "

# Row 43 - Filecoin
Set-TextValue "D43" "5.52"
Set-TextValue "E43" "  -6.05%  "

# Row 44 - Hedera
Set-TextValue "D44" "0.0857"
Set-TextValue "E44" "  -4.04%  "

# Row 45 - Mantle
Set-TextValue "D45" "0.888"
Set-TextValue "E45" "  -4.05%  "

# Row 46 - OKB
Set-TextValue "D46" "45.43"
Set-TextValue "E46" "  -3.78%  "

# Row 47 - InjectiveProtocol
Set-TextValue "D47" "27.38"
Set-TextValue "E47" "  -2.78%  "

# Row 48 - dogwifhat
Set-TextValue "D48" "2.57"
Set-TextValue "E48" "  -5.04%  "

# Row 49 - ONDO
Set-TextValue "D49" "1.21"
Set-TextValue "E49" "  -2.32%  "

# Row 50 - was SuiNetwork, now Cosmos
Set-TextValue "B50" "Cosmos"
Set-TextValue "C50" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D50" "7.52"
Set-TextValue "E50" "  -3.42%  "

# Row 51 - was Cosmos, now SuiNetwork
Set-TextValue "B51" "SuiNetwork"
Set-TextValue "C51" "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextValue "D51" "1.02"
Set-TextValue "E51" "  -4.48%  "
